# Atualização 09/07 - Programa agora verifica se a bátida atual ja foi
# realizada no mesmo dia e somente atualiza ela.
#
# The punch-clock log used to keep a separate row for every clock event of
# the day. Now, when a punch for the current day already exists, the
# existing row is simply updated instead of appending a new one - so the
# sheet shrinks from 3 data rows (4 total with header) to 2 data rows
# (3 total with header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: the 08/07/2023 entry becomes the (updated) 09/07/2023 entry ---
# Force the cell to stay text (otherwise Excel would reinterpret the
# dd/mm/yyyy-looking string as a real date value).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "09/07/2023"
$ws.Range("B2").Value = "13:32:46"
$ws.Range("C2").Value = "13:32:48"
$ws.Range("D2").Value = "13:32:49"
$ws.Range("E2").Value = "13:32:50"

# --- Row 3: becomes the new 10/07/2023 entry ---
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "10/07/2023"
$ws.Range("B3").Value = "13:35:44"
$ws.Range("C3").Value = "13:35:45"
$ws.Range("D3").Value = "13:35:46"
$ws.Range("E3").Value = "13:35:47"

# Drop the temporary text formatting again so the cells fall back to the
# workbook's default (General) style, exactly like every other cell.
$ws.Range("A2:A3").ClearFormats()

# --- Old row 4 is no longer needed; remove it entirely ---
$ws.Rows.Item(4).Delete()

Write-Host "PontoEletrônico atualizado: 09/07/2023 e 10/07/2023."
